$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking D-column cells so they keep their original text cell type
# (their new values would otherwise be auto-parsed as numbers by Excel)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range('D2').Value = '48.540.82'
$ws.Range('E2').Value = '  -2.75%  '
$ws.Range('D3').Value = '2.605.59'
$ws.Range('E3').Value = '  +1.64%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('B5').Value = 'Solana'
$ws.Range('C5').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D5').Value = '109.44'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').Value = '321.10'
$ws.Range('E6').Value = '  -0.64%  '
$ws.Range('D7').Value = '0.521'
$ws.Range('E7').Value = '  -1.55%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '0.539'
$ws.Range('E9').Value = '  -3.29%  '
$ws.Range('D10').Value = '39.20'
$ws.Range('E10').Value = '  -3.05%  '
$ws.Range('D11').Value = '19.72'
$ws.Range('E11').Value = '  -3.37%  '
$ws.Range('D12').Value = '0.0806'
$ws.Range('E12').Value = '  -1.51%  '
$ws.Range('E13').Value = '  +0.22%  '
$ws.Range('D14').Value = '7.20'
$ws.Range('E14').Value = '  -1.19%  '
$ws.Range('D15').Value = '3.009.02'
$ws.Range('E15').Value = '  +1.62%  '
$ws.Range('D16').Value = '2.599.43'
$ws.Range('E16').Value = '  +1.90%  '
$ws.Range('D17').Value = '0.860'
$ws.Range('E17').Value = '  -0.44%  '
$ws.Range('D18').Value = '48.452.25'
$ws.Range('E18').Value = '  -2.43%  '
$ws.Range('D19').Value = '2.94'
$ws.Range('E19').Value = '  -3.99%  '
$ws.Range('D20').Value = '12.79'
$ws.Range('E20').Value = '  -4.03%  '
$ws.Range('D21').Value = '6.65'
$ws.Range('E21').Value = '  -0.83%  '
$ws.Range('D22').Value = '0.0₃0941'
$ws.Range('E22').Value = '  -0.81%  '
$ws.Range('D23').Value = '269.02'
$ws.Range('E23').Value = '  -5.44%  '
$ws.Range('D24').Value = '68.48'
$ws.Range('E24').Value = '  -5.53%  '
$ws.Range('E25').Value = '  -0.58%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.22%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '25.95'
$ws.Range('E27').Value = '  -1.87%  '
$ws.Range('D28').Value = '9.98'
$ws.Range('E28').Value = '  +0.82%  '
$ws.Range('D29').Value = '2.22'
$ws.Range('E29').Value = '  -0.50%  '
$ws.Range('D30').Value = '34.74'
$ws.Range('E30').Value = '  -2.04%  '
$ws.Range('D31').Value = '0.137'
$ws.Range('E31').Value = '  -5.69%  '
$ws.Range('D32').Value = '49.23'
$ws.Range('E32').Value = '  -0.53%  '
$ws.Range('E33').Value = '  +0.41%  '
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('D35').Value = '18.94'
$ws.Range('E35').Value = '  -4.09%  '
$ws.Range('D36').Value = '0.0793'
$ws.Range('E36').Value = '  +0.67%  '
$ws.Range('D37').Value = '4.93'
$ws.Range('E37').Value = '  +4.06%  '
$ws.Range('D38').Value = '2.02'
$ws.Range('E38').Value = '  -0.76%  '
$ws.Range('D39').Value = '3.13'
$ws.Range('E39').Value = '  +3.02%  '
$ws.Range('D40').Value = '125.94'
$ws.Range('E40').Value = '  +2.28%  '
$ws.Range('E41').Value = '  -1.49%  '
$ws.Range('D42').Value = '22.03'
$ws.Range('E42').Value = '  -1.37%  '
$ws.Range('E43').Value = '  -4.19%  '
$ws.Range('D44').Value = '0.0315'
$ws.Range('E44').Value = '  +0.86%  '
$ws.Range('D45').Value = '2.055.26'
$ws.Range('E45').Value = '  +1.61%  '
$ws.Range('D46').Value = '3.22'
$ws.Range('E46').Value = '  -3.72%  '
$ws.Range('D47').Value = '2.10'
$ws.Range('E47').Value = '  +2.87%  '
$ws.Range('D48').Value = '2.17'
$ws.Range('E48').Value = '  +0.88%  '
$ws.Range('D49').Value = '8.86'
$ws.Range('E49').Value = '  -2.08%  '
$ws.Range('D50').Value = '58.37'
$ws.Range('E50').Value = '  +2.12%  '
$ws.Range('D51').Value = '5.13'
$ws.Range('E51').Value = '  -4.39%  '
